$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.981.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "'1.849.58"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.75%  "
$ws.Range("D5").Value = "'1.011"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'309.72"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "'0.4768"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("D8").Value = "'0.3673"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "'0.07215"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'0.9285"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").Value = "'19.70"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'0.07715"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "'1.830.18"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "'5.321"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "'6.421"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "'88.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "'0.000008630"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'27.022.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'14.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "'5.056"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "'1.934"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "'152.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'18.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").Value = "'1.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "'114.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'4.983"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").Value = "'0.08888"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "'3.323"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.54%  "
$ws.Range("D32").Value = "'1.172"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'0.7421"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "'4.496"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'2.722"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("D37").Value = "'0.05268"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").Value = "'0.01950"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "'2.974"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "'6.995"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").Value = "'0.1510"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "'8.196"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").Value = "'10.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.97%  "
$ws.Range("D45").Value = "'0.4733"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "'101.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("D48").Value = "'1.607"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").Value = "'65.88"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("D50").Value = "'0.06026"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'0.8870"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.08%  "
